# Add the new data row (row 8) mirroring the existing rows' layout:
#   column A -> a new ID value (new shared string)
#   column C -> "ID Used" (existing shared string, reused)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "0105052329189"
$ws.Range("C8").Value = "ID Used"

# Match the text number-format used by the other data rows in columns A/C
# (numFmtId 49 "@" — keeps the ID value stored/displayed as text).
$ws.Range("A8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
